$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(16, 1).Value = "j"
$ws.Cells.Item(16, 2).Value = 9378.0
